$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-11-30 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-01 Friday", 2) | Out-Null

# Update each answer cell in the table (row-major order, 20 rows x 5 columns)
$t = $d.Tables.Item(1)
$answers = @(
    "49-4=45",
    "98-45=53",
    "8+21=29",
    "69+10=79",
    "62+30=92",
    "2+25=27",
    "81-60=21",
    "39-20=19",
    "11+43=54",
    "58-9=49",
    "22+47=69",
    "80-56=24",
    "94-22=72",
    "83+5=88",
    "87-84=3",
    "98-14=84",
    "34+41=75",
    "60+8=68",
    "14+61=75",
    "67-51=16",
    "52+16=68",
    "75-37=38",
    "69+21=90",
    "69-15=54",
    "9+59=68",
    "5+74=79",
    "31-1=30",
    "89-33=56",
    "88-64=24",
    "34+13=47",
    "59+4=63",
    "73-67=6",
    "10+69=79",
    "92-31=61",
    "51+29=80",
    "87-10=77",
    "86-50=36",
    "58+41=99",
    "91-16=75",
    "86-18=68",
    "9+57=66",
    "13+14=27",
    "2+88=90",
    "9+83=92",
    "93-36=57",
    "56-14=42",
    "58-57=1",
    "0+91=91",
    "92-70=22",
    "42-28=14",
    "16+79=95",
    "55+11=66",
    "81-21=60",
    "93-86=7",
    "26+62=88",
    "2+38=40",
    "69+25=94",
    "1+15=16",
    "55+37=92",
    "29-5=24",
    "5+50=55",
    "36-5=31",
    "46+32=78",
    "71+24=95",
    "71-64=7",
    "35+18=53",
    "78-78=0",
    "19+49=68",
    "97-52=45",
    "70+21=91",
    "83-30=53",
    "46+47=93",
    "15+25=40",
    "94-22=72",
    "50-33=17",
    "26+10=36",
    "50-13=37",
    "61+26=87",
    "26+38=64",
    "19+7=26",
    "80-73=7",
    "58+0=58",
    "13+13=26",
    "58+32=90",
    "36+47=83",
    "85+7=92",
    "1+39=40",
    "51+8=59",
    "0+25=25",
    "25+38=63",
    "8+34=42",
    "18+9=27",
    "75-12=63",
    "95-36=59",
    "37-16=21",
    "29+17=46",
    "93-76=17",
    "44-39=5",
    "81-69=12",
    "29+24=53"
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $answers[$idx]
        $idx = $idx + 1
    }
}

Write-Output "Updated $idx cells"
